$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("20-09-2021", -11098, 2951, 918, 762, 1270),
    @("21-09-2021", -10845, 3292, 1085, 827, 1380),
    @("22-09-2021", -10932, 3769, 1163, 1239, 1367),
    @("23-09-2021", -10737, 3435, 983, 779, 1673),
    @("24-09-2021", -11108, 3289, 863, 1229, 1197),
    @("27-09-2021", -11087, 3038, 993, 962, 1083),
    @("28-09-2021", -10586, 3163, 1265, 753, 1145),
    @("29-09-2021", -10505, 3276, 912, 830, 1534),
    @("30-09-2021", -10187, 3750, 1366, 1046, 1338)
)

$startRow = 181
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
